$d = $word.ActiveDocument
try {
  Write-Output $d.CustomXMLParts.Count
} catch { Write-Output "ERR $($_.Exception.Message)" }
